# Commit: "Implement list of users on matches"
# The Software sheet gains one new row documenting the ngx-gallery npm
# package (URL hyperlink in column A, install command in column B), and
# becomes the active sheet/tab of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Software")

# New hyperlink + label for column A (reuses the existing "link" cell
# style, exactly like every other hyperlink row on this sheet).
$ws.Hyperlinks.Add($ws.Range("A11"), "https://www.npmjs.com/package/ngx-gallery")
$ws.Range("A11").Style = $ws.Range("A1").Style

# Install command for column B.
$ws.Range("B11").Value = "npm install ngx-gallery"

# Make "Software" the active sheet and put the selection on the new
# install-command cell, matching the saved workbook/sheet view state.
$ws.Activate()
$ws.Range("B11").Select()
